$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: new "RF5" column (I) next to the existing summary table (rows 46-54)
# ---------------------------------------------------------------------------
$ws.Range("I46").Value = "100-N"
$ws.Range("I47").Value = "RF5"

$ws.Range("I48").Value = 20627948
$ws.Range("I49").Value = 7
$ws.Range("I50").Value = 254352
$ws.Range("I51").Value = 2398
$ws.Range("I52").Value = 65001

# ---------------------------------------------------------------------------
# Step 2: "100-20" / "100-30" columns (J, K) plus the unlabeled L column
# ---------------------------------------------------------------------------
$ws.Range("J46").Value = "100-20"
$ws.Range("K46").Value = "100-30"

$ws.Range("J48").Value = 20646694
$ws.Range("J49").Value = 1119
$ws.Range("J50").Value = 255138
$ws.Range("J51").Value = 1831
$ws.Range("J52").Value = 65086

$ws.Range("K48").Value = 20639809
$ws.Range("K49").Value = 439
$ws.Range("K50").Value = 257051
$ws.Range("K51").Value = 1075
$ws.Range("K52").Value = 64858

$ws.Range("L48").Value = 20650981
$ws.Range("L49").Value = 23
$ws.Range("L50").Value = 256164
$ws.Range("L51").Value = 746
$ws.Range("L52").Value = 65385

# ---------------------------------------------------------------------------
# Step 3: two new confusion-matrix style blocks below (rows 59-65, 68-74)
# ---------------------------------------------------------------------------
$ws.Range("E59").Value = "yreal"

$ws.Range("A60").Value = "sql"
$ws.Range("C60").Value = 1
$ws.Range("D60").Value = 2
$ws.Range("E60").Value = 3
$ws.Range("F60").Value = 4
$ws.Range("G60").Value = 5

$ws.Range("B61").Value = 1
$ws.Range("C61").Value = 20502785
$ws.Range("D61").Value = 68436
$ws.Range("E61").Value = 11683
$ws.Range("F61").Value = 157
$ws.Range("G61").Value = 44

$ws.Range("B62").Value = 2
$ws.Range("C62").Value = 385503
$ws.Range("D62").Value = 42782
$ws.Range("E62").Value = 11223
$ws.Range("F62").Value = 121
$ws.Range("G62").Value = 56

$ws.Range("A63").Value = "ypred"
$ws.Range("B63").Value = 3
$ws.Range("C63").Value = 476660
$ws.Range("D63").Value = 167613
$ws.Range("E63").Value = 195288
$ws.Range("F63").Value = 11007
$ws.Range("G63").Value = 2861

$ws.Range("B64").Value = 4
$ws.Range("C64").Value = 11492
$ws.Range("D64").Value = 7467
$ws.Range("E64").Value = 50920
$ws.Range("F64").Value = 17726
$ws.Range("G64").Value = 5938

$ws.Range("B65").Value = 5
$ws.Range("C65").Value = 9221
$ws.Range("D65").Value = 3024
$ws.Range("E65").Value = 30824
$ws.Range("F65").Value = 31482
$ws.Range("G65").Value = 65890

$ws.Range("E68").Value = "yreal"

$ws.Range("A69").Value = "RF"
$ws.Range("C69").Value = 1
$ws.Range("D69").Value = 2
$ws.Range("E69").Value = 3
$ws.Range("F69").Value = 4
$ws.Range("G69").Value = 5

$ws.Range("B70").Value = 1
$ws.Range("C70").Value = 20650981
$ws.Range("D70").Value = 81531
$ws.Range("E70").Value = 14611
$ws.Range("F70").Value = 190
$ws.Range("G70").Value = 56

$ws.Range("B71").Value = 2
$ws.Range("C71").Value = 202
$ws.Range("D71").Value = 23
$ws.Range("E71").Value = 1
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0

$ws.Range("A72").Value = "ypred"
$ws.Range("B72").Value = 3
$ws.Range("C72").Value = 725433
$ws.Range("D72").Value = 204885
$ws.Range("E72").Value = 256164
$ws.Range("F72").Value = 29429
$ws.Range("G72").Value = 9046

$ws.Range("B73").Value = 4
$ws.Range("C73").Value = 131
$ws.Range("D73").Value = 48
$ws.Range("E73").Value = 813
$ws.Range("F73").Value = 746
$ws.Range("G73").Value = 302

$ws.Range("B74").Value = 5
$ws.Range("C74").Value = 8914
$ws.Range("D74").Value = 2835
$ws.Range("E74").Value = 28349
$ws.Range("F74").Value = 30128
$ws.Range("G74").Value = 65385

# ---------------------------------------------------------------------------
# Step 4: "combine1" / "combine2" / "combine3" columns (M, N, O)
# ---------------------------------------------------------------------------
$ws.Range("M46").Value = "combine1"
$ws.Range("N46").Value = "combine2"
$ws.Range("O46").Value = "combine3"

$ws.Range("M48").Value = 20502787
$ws.Range("M49").Value = 42782
$ws.Range("M50").Value = 195288
$ws.Range("M51").Value = 17576
$ws.Range("M52").Value = 66052

$ws.Range("N48").Value = 20502781
$ws.Range("N49").Value = 13118
$ws.Range("N50").Value = 254329
$ws.Range("N51").Value = 97
$ws.Range("N52").Value = 66052

# ---------------------------------------------------------------------------
# Row 53: re-enter the aggregate formula across the whole A:N range so it
# becomes one shared formula (matches the rest of the already-existing row).
# ---------------------------------------------------------------------------
$ws.Range("A53:N53").Formula = "=A48+A49*10+A50*50+A51*100+A52*200"

# Row 54: ratio formulas for the new I..N columns (same pattern as B54:H54)
$ws.Range("I54").Formula = '=I53/$A$53'
$ws.Range("J54").Formula = '=J53/$A$53'
$ws.Range("K54").Formula = '=K53/$A$53'
$ws.Range("L54").Formula = '=L53/$A$53'
$ws.Range("M54").Formula = '=M53/$A$53'
$ws.Range("N54").Formula = '=N53/$A$53'

# ---------------------------------------------------------------------------
# Update the view to match the final saved selection/scroll position
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("O47").Select()
